$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header cells: "_old" -> "_FV2410", "_new" -> "_FV2504"
$fv2410Headers = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")
$fv2504Headers = @("Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504","Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410Headers[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504Headers[$i]
}

# 2) Turn the used range into an Excel Table ("Table1") spanning A1:U65
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 3) Freeze the header row (split after row 1, freeze)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
